$d = $word.ActiveDocument

# Helper: replace a whole-word number token inside a specific paragraph,
# verifying the paragraph's current text first so we fail loudly instead
# of silently touching the wrong paragraph if numbering ever shifts.
function Update-Count($paraIndex, $expectedSubstring, $oldVal, $newVal) {
    $p = $d.Paragraphs($paraIndex)
    $text = $p.Range.Text
    if ($text -notlike "*$expectedSubstring*") {
        Write-Host "WARNING: paragraph $paraIndex does not contain '$expectedSubstring' (got: $text)"
    }
    $ok = $p.Range.Find.Execute($oldVal, $true, $true, $false, $false, $false, $true, 1, $false, $newVal, 2)
    Write-Host "Paragraph $paraIndex : '$oldVal' -> '$newVal' => $ok (now: $($p.Range.Text))"
}

# Section "фактическое количество машин:"
# "168ч. – 4;" -> "168ч. – 3;"
Update-Count 10 "168" "4" "3"

# "79ч. – 4;" -> "79ч. – 6;"
Update-Count 11 "79" "4" "6"

# Section "фактическое количество нехватки машин:"
# "168ч. – 1;" -> "168ч. – 0;"
Update-Count 18 "168" "1" "0"

# Section "планируемое количество нехватки машин:"
# "168ч. – 1;" -> "168ч. – 0;"
Update-Count 34 "168" "1" "0"
